$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Receptor-side / Edge values updated with new TPM-derived figures)
$ws.Range("M2").Value = 10.34761366666667
$ws.Range("N2").Value = 31.042841
$ws.Range("O2").Value = 0.2299953477621856
$ws.Range("P2").Value = 0.2299953477621856
$ws.Range("Q2").Value = 13.72675108088067
$ws.Range("R2").Value = 123.540759727926
$ws.Range("S2").Value = 0.2299953477621856
$ws.Range("T2").Value = 0.2299953477621856

# Row 3 (only specificity columns change)
$ws.Range("O3").Value = 0.6794731949692173
$ws.Range("P3").Value = 0.6794731949692174
$ws.Range("S3").Value = 0.6794731949692173
$ws.Range("T3").Value = 0.6794731949692174

# Row 4
$ws.Range("M4").Value = 4.073058666666666
$ws.Range("N4").Value = 12.219176
$ws.Range("O4").Value = 0.09053145726859702
$ws.Range("P4").Value = 0.09053145726859703
$ws.Range("Q4").Value = 5.403164850970667
$ws.Range("R4").Value = 48.62848365873599
$ws.Range("S4").Value = 0.09053145726859702
$ws.Range("T4").Value = 0.09053145726859703

$wb.Save()
